# Fruta / hortaliza, semanal
# Insert a new weekly record at row 114 ("Cebollín" dataset), pushing the
# existing rows 114:133 down to 115:134.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 114:133 down by one row (creates a fresh, empty row 114 and
# extends the used range to row 134).
$ws.Rows("114:114").Insert()

# Populate the newly inserted row 114 with this week's data point.
$ws.Range("A114").Value = 8
$ws.Range("B114").Value = "Terminal La Palmera de La Serena"
$ws.Range("C114").Value = "Coquimbo"
$ws.Range("D114").Value = 44522
$ws.Range("E114").Value = 4
$ws.Range("F114").Value = 100112037
$ws.Range("G114").Value = "Cebollín"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 2900
$ws.Range("K114").Value = 900
$ws.Range("L114").Value = 1000
$ws.Range("M114").Value = 950
$ws.Range("N114").Value = "$/paquete 6 unidades"
$ws.Range("O114").Value = "Provincia del Elquí"
$ws.Range("P114").Value = 158
$ws.Range("Q114").Value = 6
$ws.Range("R114").Value = "Hortaliza"

# Make sure the date column keeps the same date number format as the rest
# of the column.
$ws.Range("D114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
